$d = $word.ActiveDocument

# --- Locate the anchor paragraph: the one containing the ellipsis "..." -
$ellipsisChar = [char]0x2026
$ellipsisText = "$ellipsisChar`r"

$ellipsisIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -eq $ellipsisText) {
        $ellipsisIndex = $i
    }
}
if ($ellipsisIndex -eq -1) {
    # Fallback: the ellipsis paragraph is the second-to-last paragraph.
    $ellipsisIndex = $d.Paragraphs.Count - 1
}

# --- 1) Add the "To complete the assignment..." paragraph ---------------
$pEllipsis = $d.Paragraphs.Item($ellipsisIndex)
$pEllipsis.Range.InsertParagraphAfter()

$pConfirm = $d.Paragraphs.Item($ellipsisIndex + 1)
$r = $pConfirm.Range
$r.Collapse(1)
$r.InsertAfter("To")
$r.Collapse(0)
$r.InsertAfter(" complete the assignment, I will make changes in this document to prove that I have completed the assignment.")
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("I will write my username to my submission document.")

# --- 2) Add the "Yu Zhao" heading paragraph after the trailing blank ----
# The trailing empty paragraph (already present at the end of the body)
# stays untouched; the new heading paragraph goes right after it.
$blankIndex = $ellipsisIndex + 2
$pBlank = $d.Paragraphs.Item($blankIndex)
$pBlank.Range.InsertParagraphAfter()

$pName = $d.Paragraphs.Item($blankIndex + 1)
$pName.Range.Text = "Yu Zhao"
$pName.Style = "Heading 1"
$pName.Range.LanguageIDFarEast = "zh-CN"

Write-Output "Edit complete"
